$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time costs")

# Excel alignment constants
$xlLeft = -4131
$xlCenter = -4108

# New row 9 (row 8 is intentionally left blank, matching the source data)
$ws.Range("A9").Value = 41689
$ws.Range("A9").NumberFormat = "dd/mm/yy;@"
$ws.Range("A9").HorizontalAlignment = $xlLeft

$ws.Range("B9").Value = 4
$ws.Range("B9").HorizontalAlignment = $xlLeft

$ws.Range("C9").Value = "x"
$ws.Range("C9").HorizontalAlignment = $xlCenter
$ws.Range("C9").VerticalAlignment = $xlCenter

$ws.Range("D9").Value = "x"
$ws.Range("D9").HorizontalAlignment = $xlCenter
$ws.Range("D9").VerticalAlignment = $xlCenter

$ws.Range("E9").Value = "x"
$ws.Range("E9").HorizontalAlignment = $xlCenter
$ws.Range("E9").VerticalAlignment = $xlCenter

$ws.Range("F9").Value = "x"
$ws.Range("F9").HorizontalAlignment = $xlCenter
$ws.Range("F9").VerticalAlignment = $xlCenter

$ws.Range("G9").Value = "GUI development"
$ws.Range("G9").HorizontalAlignment = $xlLeft
$ws.Range("G9").VerticalAlignment = $xlCenter

# Update the active selection to match the author's final cursor position
$ws.Range("G11").Select()
